# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header values -------------------------------------------------
# VALOR MORA total (sum of Valor Mora column for the new, shorter table)
$ws.Range("E11").Value = 113552
# Cant. Periodos (distinct periods now represented: 2009, 2010, 2012)
$ws.Range("F13").Value = 3

# --- Replace the detail table (rows 16-27) with the new, shorter dataset ---
# Row 16: Robinson Felipe Gelvis Pacheco - periodo 2009
$ws.Range("C16").Value = "1143359583"
$ws.Range("D16").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E16").Value = "2009"
$ws.Range("F16").Value = 38000
$ws.Range("G16").Value = 950000

# Row 17: Robinson Felipe Gelvis Pacheco - periodo 2010
$ws.Range("C17").Value = "1143359583"
$ws.Range("D17").Value = "ROBINSON FELIPE GELVIS PACHECO"
$ws.Range("E17").Value = "2010"
$ws.Range("F17").Value = 38000
$ws.Range("G17").Value = 950000

# Row 18: Asmeth Leonar Marrugo Gonzalez - periodo 2012
$ws.Range("C18").Value = "73350807"
$ws.Range("D18").Value = "ASMETH LEONAR MARRUGO GONZALEZ"
$ws.Range("E18").Value = "2012"
$ws.Range("F18").Value = 1211
$ws.Range("G18").Value = 908526

# Row 19: Luis Tomas Causil Paternina - periodo 2012 (new last row of the table)
$ws.Range("C19").Value = "1051444928"
$ws.Range("D19").Value = "LUIS TOMAS CAUSIL PATERNINA"
$ws.Range("E19").Value = "2012"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

# Row 19 becomes the bottom-most row of the table, so it needs the "closing"
# border formatting that row 27 (the previous bottom row) carries. Copy just
# the formatting before the now-unused rows 20-27 are removed.
$ws.Range("B27:J27").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Remove the now-obsolete extra rows (20-27); this shifts the footer
# ("___...", "NOMBRE DEL REPRESENTANTE LEGAL", "FIRMA DEL REPRESENTANTE LEGAL")
# up from rows 32/33 to rows 24/25.
$ws.Rows("20:27").Delete()
